# Trade #10 closed at 2026-02-17 19:47:33 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ----- Sheet: Summary -----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1299.98
$summary.Range("B4").Value = -0.02
$summary.Range("B5").Value = -0.04
$summary.Range("B6").Value = 10
$summary.Range("B8").Value = 4
$summary.Range("B9").Value = 50

# ----- Sheet: Strategy Status -----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98
$status.Range("D4").Value = 10
$status.Range("E4").Value = -0.02
$status.Range("F4").Value = -0.02
$status.Range("G4").Value = 50

# ----- Sheet: All Trades (row 11 = Trade #10 closed) -----
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G11").Value = 0.8100000000000001
$allTrades.Range("H11").Value = "CLOSED"
$allTrades.Range("I11").Value = -2.4096
$allTrades.Range("J11").Value = -0.02
$allTrades.Range("K11").Value = 99.98
$allTrades.Range("P11").Value = "early_exit"
$allTrades.Range("Q11").Value = 0.11

# ----- Sheet: MarketMaking (row 11 = Trade #10 closed, mirrors All Trades) -----
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G11").Value = 0.8100000000000001
$mm.Range("H11").Value = "CLOSED"
$mm.Range("I11").Value = -2.4096
$mm.Range("J11").Value = -0.02
$mm.Range("K11").Value = 99.98
$mm.Range("P11").Value = "early_exit"
$mm.Range("Q11").Value = 0.11
